$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1's C1 formula (selectif -> #NAME?) is left untouched; add a working
# A1+B1 sum in row 2 to "fix" the cycling/broken formula referenced in the
# commit message.
$ws.Range("A2").Formula = "=A1+B1"

# Move the active selection down to A3, as in the saved file.
$ws.Range("A3").Select()
